# 12/02 - New Branch commit
# Applies the testData.xlsx edits:
#  - Environment_DirectSales (sheet1): B3 gets the Hyperlink style + a
#    mailto hyperlink, B4's value moves from "Testing-2021" to
#    "Testing-2022", and the sheet becomes the active/selected tab with
#    B5 selected.
#  - GeneralVariables (sheet3): a new company row (SOI79) is inserted
#    right after the existing SOI69 row, three brand-new rows are
#    appended at the bottom, and the sheet's own selection moves to F13
#    (losing the "active tab" flag, since DirectSales becomes active).

$wb = $excel.ActiveWorkbook

$wsDirect = $wb.Worksheets.Item("Environment_DirectSales")
$wsGeneral = $wb.Worksheets.Item("GeneralVariables")

# --- Environment_DirectSales --------------------------------------------

# B4: Testing-2021 -> Testing-2022
$wsDirect.Range("B4").Value = "Testing-2022"

# B3: add a mailto hyperlink for the e-mail address already stored there,
# then re-apply the Hyperlink cell style (Add() resets the style id).
$wsDirect.Hyperlinks.Add($wsDirect.Range("B3"), "mailto:uft.test.automation@gmail.com")
$wsDirect.Range("B3").Style = "Hyperlink"

# --- GeneralVariables -----------------------------------------------------

# Insert a brand-new "SOI79" company-name row right after row 5 (the SOI69
# name row), pushing the id block (old rows 6-9) and everything else down
# by one.
$wsGeneral.Rows("6:6").Insert()
$wsGeneral.Range("A6").Value = "testingCompanySOI79"
$wsGeneral.Range("B6").Value = "AutoTestingCompany_SOI79"

# Insert the matching "SOI79" id row right after the (now shifted) id block,
# which ends at row 10 (old row 9) -> new row 11.
$wsGeneral.Rows("11:11").Insert()
$wsGeneral.Range("A11").Value = "idTestingCompanySOI79"
$wsGeneral.Range("B11").Value = "0013E00000zZ4ccQAC"

# Append the new rows at the bottom of the sheet (old last row 18 is now 20
# after the two inserts above).
$wsGeneral.Range("A21").Value = "noServiceErrorMessage"
$wsGeneral.Range("B21").Value = "This order can't be submited!No services associated to the order."

$wsGeneral.Range("A22").Value = "noFilesOnServiceErrorMessage"
$wsGeneral.Range("B22").Value = "This order can't be submitted.One or more services without files"

$wsGeneral.Range("A23").Value = "companyContactPerson"
$wsGeneral.Range("B23").Value = "André Esteves"

# --- Selections / active sheet --------------------------------------------

# GeneralVariables keeps workbookViewId=0 but its selection moves to F13
# (and it stops being the tab-selected sheet once DirectSales is activated).
$wsGeneral.Range("F13").Select()

# Environment_DirectSales becomes the active/selected tab with B5 selected.
$wsDirect.Activate()
$wsDirect.Range("B5").Select()
